$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 2 and row 5 (columns A, C, D, E, F — B is unchanged in both rows)
$ws.Range("A2").Value = "a03_01"
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 19

$ws.Range("A5").Value = "a03_04"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 15

# Update the selection to match the saved state
$ws.Range("A2:XFD4").Select()
$ws.Cells.Item(4, 1).Activate()
